# Applies the CMDPII notification corrections described in the diff:
# student name, class, dates, infraction article, disciplinary note
# deltas/classification, and the closing date line.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "Aluno: chiquin" "Aluno: Quirino Alessandro Cordeiro Gomes"
Replace-Text "Turma: 6ºA" "Turma: 2ºA"
Replace-Text "Data: 23 de junho de 2025" "Data: 15 de junho de 2025"
Replace-Text "Art. 54 – 14" "Art. 54 – 13"
Replace-Text "Esta medida acarreta perda de sua nota disciplinar em -0.50 pontos, enquadrando-se no comportamento Excepcional." "Esta medida acarreta perda de sua nota disciplinar em -2.10 pontos, enquadrando-se no comportamento Insuficiente."
Replace-Text "Nota Anterior: 10.00" "Nota Anterior: 6.10"
Replace-Text "Nota Atual: 9.50" "Nota Atual: 4.00"
Replace-Text "Cruzeiro do Sul – AC, 23 de junho de 2025" "Cruzeiro do Sul – AC, 15 de junho de 2025"
